$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The roster grew (new students) and was re-sorted alphabetically by name.
# Write the brand-new names first, in the precise order they were first
# introduced, so the workbook shared-string table builds up in the same
# append order as the authored edit.
$ws.Cells.Item(19, 2).Value = "Silvia Naghi"
$ws.Cells.Item(17, 2).Value = "Raul Andrei"
$ws.Cells.Item(12, 2).Value = "Levente Nagy"
$ws.Cells.Item(14, 2).Value = "Mark Pop"
$ws.Cells.Item(3, 2).Value = "Alessandro Vereș-Pop"
$ws.Cells.Item(9, 2).Value = "Daniela Cionca (Mărie)"
$ws.Cells.Item(13, 2).Value = "Luca Șeicaru"
$ws.Cells.Item(6, 2).Value = "Attila Bunta"
$ws.Cells.Item(8, 2).Value = "Codruț Avram"
$ws.Cells.Item(21, 2).Value = "Victor Lazăr"

# Now fill in the rest of the roster (names already present earlier in the
# sheet, just relocated to their new alphabetically-sorted row).
$ws.Cells.Item(4, 2).Value = "Amanda Hajdu"
$ws.Cells.Item(5, 2).Value = "Andra Agud"
$ws.Cells.Item(7, 2).Value = "Claudiu Druța"
$ws.Cells.Item(10, 2).Value = "Delia Negrea"
$ws.Cells.Item(11, 2).Value = "Denisa Cioban"
$ws.Cells.Item(15, 2).Value = "Miriam Bacso"
$ws.Cells.Item(16, 2).Value = "Paul Dobroțchi"
$ws.Cells.Item(18, 2).Value = "Răzvan Baroi"
$ws.Cells.Item(20, 2).Value = "Sorin Fechete"

# Attendance ("Prezente") counts for each student in column C.
$ws.Cells.Item(3, 3).Value = 2
$ws.Cells.Item(4, 3).Value = 1
$ws.Cells.Item(5, 3).Value = 1
$ws.Cells.Item(6, 3).Value = 1
$ws.Cells.Item(7, 3).Value = 2
$ws.Cells.Item(8, 3).Value = 1
$ws.Cells.Item(9, 3).Value = 2
$ws.Cells.Item(10, 3).Value = 2
$ws.Cells.Item(11, 3).Value = 1
$ws.Cells.Item(12, 3).Value = 2
$ws.Cells.Item(13, 3).Value = 1
$ws.Cells.Item(14, 3).Value = 2
$ws.Cells.Item(15, 3).Value = 2
$ws.Cells.Item(16, 3).Value = 2
$ws.Cells.Item(17, 3).Value = 1
$ws.Cells.Item(18, 3).Value = 2
$ws.Cells.Item(19, 3).Value = 2
$ws.Cells.Item(20, 3).Value = 1
$ws.Cells.Item(21, 3).Value = 1

# Restore the active-cell selection left by the editor.
$ws.Range("H12").Select() | Out-Null
